# Generate Report for Archive
#
# The localization run moved on from "Ready for handoff" to "In Translation",
# so the Status cells that held the old text need updating everywhere they
# appear (the Overview rollup sheet plus each per-locale detail sheet), and
# the now-narrower "Status" column needs to shrink to match.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" -----------------
# Overview!E2 is the zh-cn status, Overview!F2 is the de-de status.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Each locale detail sheet carries its own Status cell in column C.
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the Status column(s) -------------------------------------------
# Overview columns E (zh-cn) and F (de-de), and column C ("Status") on each
# locale detail sheet. ColumnWidth is in characters; Excel stores/quantizes
# the underlying width to whole pixels, so 12.5 characters is the closest
# achievable width to the target ~13.41-character column width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
